$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename column headers to use the respective input-file-name suffixes
#    ("_old"/"_new" -> "_FV2310"/"_FV2404") instead of the generic ones.
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J => "_FV2310" suffix, column K is the "diff" column (unchanged),
# columns L-U => "_FV2404" suffix.
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $leftCol = [char]([int][char]'A' + $i)
    $ws.Range("$leftCol`1").Value = "$($baseNames[$i])_FV2310"

    $rightCol = [char]([int][char]'L' + $i)
    $ws.Range("$rightCol`1").Value = "$($baseNames[$i])_FV2404"
}

# 2. Turn the data range into an actual Excel Table, using the (renamed)
#    header row as the column names.
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U59"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# 3. Freeze the header row so it stays visible while scrolling.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
